$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 4123.76
$ws.Range("I33").Value = 5989.0586
$ws.Range("K33").Value = 5989.0586
$ws.Range("M33").Value = -5760.0586

$ws.Range("H86").Value = 28046.908
$ws.Range("I86").Value = 16153.385
$ws.Range("J86").Value = 45226.445
$ws.Range("K86").Value = 16153.385
$ws.Range("L86").Value = 45226.445
$ws.Range("M86").Value = -15030.385
$ws.Range("N86").Value = -47472.445

$ws.Range("H89").Value = 28046.908
$ws.Range("I89").Value = 16153.385
$ws.Range("J89").Value = 45226.445
$ws.Range("K89").Value = 80766.925
$ws.Range("L89").Value = 226132.225
$ws.Range("M89").Value = -75150.925
$ws.Range("N89").Value = -237364.225

$ws.Range("H98").Value = 2025.5
$ws.Range("I98").Value = 2516.6667
$ws.Range("J98").Value = 1534.3334
$ws.Range("K98").Value = 2516.6667
$ws.Range("L98").Value = 1534.3334
$ws.Range("M98").Value = -1018.6667
$ws.Range("N98").Value = -4530.3334

$ws.Range("H100").Value = 40904.42
$ws.Range("I100").Value = 2022.3889
$ws.Range("J100").Value = 128389
$ws.Range("K100").Value = 2022.3889
$ws.Range("L100").Value = 128389
$ws.Range("M100").Value = -1481.3889
$ws.Range("N100").Value = -129471

$ws.Range("H112").Value = 1144.4286
$ws.Range("J112").Value = 1150.3903
$ws.Range("L112").Value = 3451.1709
$ws.Range("N112").Value = -5667.1709

$ws.Range("H122").Value = 2025.5
$ws.Range("I122").Value = 2516.6667
$ws.Range("J122").Value = 1534.3334
$ws.Range("K122").Value = 7550.000100000001
$ws.Range("L122").Value = 4603.0002
$ws.Range("M122").Value = -5100.000100000001
$ws.Range("N122").Value = -9503.0002

$ws.Range("H129").Value = 19425824
$ws.Range("I129").Value = 494.83334
$ws.Range("J129").Value = 21959562
$ws.Range("K129").Value = 1484.50002
$ws.Range("L129").Value = 65878686
$ws.Range("M129").Value = 3515.49998
$ws.Range("N129").Value = -65888686

$ws.Range("H132").Value = 36104.465
$ws.Range("I132").Value = 42793.36
$ws.Range("J132").Value = 2660
$ws.Range("K132").Value = 128380.08
$ws.Range("L132").Value = 7980
$ws.Range("M132").Value = -125850.08
$ws.Range("N132").Value = -13040

$ws.Range("H137").Value = 885.3
$ws.Range("J137").Value = 890.9091
$ws.Range("L137").Value = 2672.7273
$ws.Range("N137").Value = -7772.7273

$ws.Range("H138").Value = 2961.65
$ws.Range("I138").Value = 1428.7059
$ws.Range("J138").Value = 4557.163
$ws.Range("K138").Value = 4286.1177
$ws.Range("L138").Value = 13671.489
$ws.Range("M138").Value = 853.8823000000002
$ws.Range("N138").Value = -23951.489

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1491
$ws.Range("I2").Value = 1486.55
$ws.Range("J2").Value = 1502.125
$ws.Range("K2").Value = 1486.55
$ws.Range("L2").Value = 1502.125
$ws.Range("M2").Value = -1373.55
$ws.Range("N2").Value = -1728.125

$ws.Range("H32").Value = 6051.5557
$ws.Range("I32").Value = 5348.1816
$ws.Range("J32").Value = 37000
$ws.Range("K32").Value = 5348.1816
$ws.Range("L32").Value = 37000
$ws.Range("M32").Value = -5061.1816
$ws.Range("N32").Value = -37574

$ws.Range("H74").Value = 1427.2972
$ws.Range("I74").Value = 970.0833
$ws.Range("J74").Value = 2271.3845
$ws.Range("K74").Value = 970.0833
$ws.Range("L74").Value = 2271.3845
$ws.Range("M74").Value = -96.08330000000001
$ws.Range("N74").Value = -4019.3845

$ws.Range("H77").Value = 1427.2972
$ws.Range("I77").Value = 970.0833
$ws.Range("J77").Value = 2271.3845
$ws.Range("K77").Value = 4850.4165
$ws.Range("L77").Value = 11356.9225
$ws.Range("M77").Value = -482.4165000000003
$ws.Range("N77").Value = -20092.9225

$ws.Range("H110").Value = 87087
$ws.Range("I110").Value = 100504.4
$ws.Range("K110").Value = 100504.4
$ws.Range("M110").Value = -98459.39999999999

$ws.Range("H116").Value = 1491
$ws.Range("I116").Value = 1486.55
$ws.Range("J116").Value = 1502.125
$ws.Range("K116").Value = 1486.55
$ws.Range("L116").Value = 1502.125
$ws.Range("M116").Value = 807.45
$ws.Range("N116").Value = -6090.125

$ws.Range("H122").Value = 1875.579
$ws.Range("I122").Value = 1643.7241
$ws.Range("J122").Value = 2622.6667
$ws.Range("K122").Value = 4931.1723
$ws.Range("L122").Value = 7868.000100000001
$ws.Range("M122").Value = -2481.1723
$ws.Range("N122").Value = -12768.0001

$ws.Range("H132").Value = 2315.4546
$ws.Range("I132").Value = 1657.1111
$ws.Range("J132").Value = 2771.2307
$ws.Range("K132").Value = 4971.3333
$ws.Range("L132").Value = 8313.6921
$ws.Range("M132").Value = -2441.3333
$ws.Range("N132").Value = -13373.6921

$ws.Range("H139").Value = 45516.25
$ws.Range("J139").Value = 45516.25
$ws.Range("L139").Value = 45516.25
$ws.Range("N139").Value = -55796.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1491
$ws.Range("I3").Value = 1486.55
$ws.Range("J3").Value = 1502.125
$ws.Range("K3").Value = 1486.55
$ws.Range("L3").Value = 1502.125
$ws.Range("M3").Value = -1372.55
$ws.Range("N3").Value = -1730.125

$ws.Range("H112").Value = 0
$ws.Range("J112").Value = 0
$ws.Range("L112").Value = 0
$ws.Range("N112").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 10832.909
$ws.Range("I99").Value = 1632.75
$ws.Range("J99").Value = 35366.668
$ws.Range("K99").Value = 1632.75
$ws.Range("L99").Value = 35366.668
$ws.Range("M99").Value = -134.75
$ws.Range("N99").Value = -38362.668

$ws.Range("H122").Value = 772243.5600000001
$ws.Range("I122").Value = 1924.8
$ws.Range("J122").Value = 1253692.8
$ws.Range("K122").Value = 5774.4
$ws.Range("L122").Value = 3761078.4
$ws.Range("M122").Value = -3324.4
$ws.Range("N122").Value = -3765978.4

$ws.Range("H126").Value = 10832.909
$ws.Range("I126").Value = 1632.75
$ws.Range("J126").Value = 35366.668
$ws.Range("K126").Value = 4898.25
$ws.Range("L126").Value = 106100.004
$ws.Range("M126").Value = -2428.25
$ws.Range("N126").Value = -111040.004

$ws.Range("H132").Value = 2679.5588
$ws.Range("I132").Value = 2081.5
$ws.Range("J132").Value = 3776
$ws.Range("K132").Value = 6244.5
$ws.Range("L132").Value = 11328
$ws.Range("M132").Value = -3714.5
$ws.Range("N132").Value = -16388

$ws.Range("H134").Value = 3187.7144
$ws.Range("I134").Value = 4071.4
$ws.Range("J134").Value = 978.5
$ws.Range("K134").Value = 12214.2
$ws.Range("L134").Value = 2935.5
$ws.Range("M134").Value = -9679.200000000001
$ws.Range("N134").Value = -8005.5

$ws.Range("H140").Value = 89672.664
$ws.Range("J140").Value = 89672.664
$ws.Range("L140").Value = 89672.664
$ws.Range("N140").Value = -100032.664

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 1346.125
$ws.Range("I122").Value = 490.85715
$ws.Range("J122").Value = 7333
$ws.Range("K122").Value = 4417.71435
$ws.Range("L122").Value = 65997
$ws.Range("M122").Value = -1967.71435
$ws.Range("N122").Value = -70897

$ws.Range("H127").Value = 955.2857
$ws.Range("J127").Value = 955.2857
$ws.Range("L127").Value = 2865.8571
$ws.Range("N127").Value = -12785.8571

$ws.Range("H132").Value = 767.5
$ws.Range("I132").Value = 593.9286
$ws.Range("J132").Value = 1375
$ws.Range("K132").Value = 5345.3574
$ws.Range("L132").Value = 12375
$ws.Range("M132").Value = -2815.3574
$ws.Range("N132").Value = -17435

$ws.Range("H133").Value = 1981.579
$ws.Range("I133").Value = 1961.1111
$ws.Range("K133").Value = 5883.3333
$ws.Range("M133").Value = -823.3333000000002

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H133").Value = 45000
$ws.Range("J133").Value = 45000
$ws.Range("L133").Value = 45000
$ws.Range("N133").Value = -55120

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 39845.85
$ws.Range("I7").Value = 57904.945
$ws.Range("J7").Value = 3727.6667
$ws.Range("K7").Value = 57904.945
$ws.Range("L7").Value = 3727.6667
$ws.Range("M7").Value = -57792.945
$ws.Range("N7").Value = -3951.6667

$ws.Range("H95").Value = 0
$ws.Range("J95").Value = 0
$ws.Range("L95").Value = 0
$ws.Range("N95").ClearContents()

$ws.Range("H126").Value = 39845.85
$ws.Range("I126").Value = 57904.945
$ws.Range("J126").Value = 3727.6667
$ws.Range("K126").Value = 173714.835
$ws.Range("L126").Value = 11183.0001
$ws.Range("M126").Value = -171244.835
$ws.Range("N126").Value = -16123.0001

$ws.Range("H132").Value = 10542.889
$ws.Range("I132").Value = 8940.258
$ws.Range("J132").Value = 14091.571
$ws.Range("K132").Value = 26820.774
$ws.Range("L132").Value = 42274.713
$ws.Range("M132").Value = -24290.774
$ws.Range("N132").Value = -47334.713

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 31671.428
$ws.Range("J54").Value = 31671.428
$ws.Range("L54").Value = 31671.428
$ws.Range("N54").Value = -32711.428

$ws.Range("H97").Value = 0
$ws.Range("J97").Value = 0
$ws.Range("L97").Value = 0
$ws.Range("N97").ClearContents()

$ws.Range("H122").Value = 69186
$ws.Range("I122").Value = 85265
$ws.Range("J122").Value = 4870
$ws.Range("K122").Value = 255795
$ws.Range("L122").Value = 14610
$ws.Range("M122").Value = -253345
$ws.Range("N122").Value = -19510
